# Fix issue with sample numbers 001-006 causing false positive data columns
# to be included in PCA plots. Update the "UnitMass" (column C) values on
# Sheet1 for both the "+ loading" and "- loading" tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 59
    3  = 73
    4  = 102
    5  = 83
    6  = 175
    7  = 46
    8  = 31
    9  = 74
    10 = 125
    11 = 103
    12 = 71
    13 = 159
    14 = 97
    15 = 58
    16 = 18
    17 = 115
    18 = 231
    19 = 72
    20 = 119
    21 = 60
    23 = 27
    25 = 39
    26 = 43
    27 = 67
    28 = 123
    29 = 79
    30 = 53
    31 = 24
    32 = 91
    33 = 121
    34 = 93
    35 = 65
    36 = 153
    37 = 94
    38 = 28
    39 = 57
    40 = 80
    41 = 52
    42 = 66
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}
